$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "doc_ids" column (E) was mostly placeholder-empty (wrapped style with no
# value). Clear the whole column's contents+formatting first so the leftover
# empty cells disappear entirely, then re-populate the rows that now carry
# real tag/document-id values (new demo data).
$ws.Range("E2:E43").Clear()

$ws.Range("E6").Value = "statpop-info"
$ws.Range("E7").Value = "pdf-wiki"
$ws.Range("E8").Value = "pdf_online, bevnat-variable"
$ws.Range("E12").Value = "pdf-wiki, tourisme-exemple"
$ws.Range("E14").Value = "bevnat-variable"
$ws.Range("E18").Value = "statpop-info, tourisme-exemple"
$ws.Range("E19").Value = "pop-com-1, pdf_online"
$ws.Range("E22").Value = "pdf-wiki"
$ws.Range("E37").Value = "tourisme-exemple"
$ws.Range("E38").Value = "bevnat-info, tourisme-exemple"
$ws.Range("E40").Value = "statpop-info"

# Restore the view to the top of the sheet, with D4 selected (instead of the
# previously scrolled-down / bottom-of-sheet selection).
$ws.Range("D4").Select()
